$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the new "B3" value first so that the new shared strings are
# interned in the same order Excel produced them in the target workbook:
# hi, userName, password, buy (existing strings Sankar/Naman keep their
# original relative order and are simply reused).
$ws.Range("B3").Value = "hi"
$ws.Range("A1").Value = "userName"
$ws.Range("B1").Value = "password"
$ws.Range("B2").Value = "buy"
$ws.Range("A2").Value = "Sankar"
$ws.Range("A3").Value = "Naman"

# Widen column B to fit the new "password"/"userName" header text and
# select B2 as the active cell (matching the authored workbook state).
$ws.Columns.Item(2).ColumnWidth = 25

$ws.Range("B2").Select()
